# Applies the scheduled-runner price/profit refresh to the Leve profit sheets.
# Each sheet (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) holds per-leve market-board
# derived columns H:N (currentAveragePrice.., LevePrice.., LeveProfit..) that are
# refreshed from an external price feed on a schedule; this updates the handful of
# rows whose cached figures changed.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 4157.448
$ws.Range("I64").Value = 3700.2222
$ws.Range("J64").Value = 4363.2
$ws.Range("K64").Value = 3700.2222
$ws.Range("L64").Value = 4363.2
$ws.Range("M64").Value = -3452.2222
$ws.Range("N64").Value = -4859.2
$ws.Range("H67").Value = 4157.448
$ws.Range("I67").Value = 3700.2222
$ws.Range("J67").Value = 4363.2
$ws.Range("K67").Value = 3700.2222
$ws.Range("L67").Value = 4363.2
$ws.Range("M67").Value = -2842.2222
$ws.Range("N67").Value = -6079.2
$ws.Range("H137").Value = 4001908
$ws.Range("I137").Value = 9092133
$ws.Range("J137").Value = 2445.8572
$ws.Range("K137").Value = 27276399
$ws.Range("L137").Value = 7337.571599999999
$ws.Range("M137").Value = -27273849
$ws.Range("N137").Value = -12437.5716

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 125263450
$ws.Range("J61").Value = 550007
$ws.Range("L61").Value = 550007
$ws.Range("N61").Value = -550431
$ws.Range("H74").Value = 12601015
$ws.Range("I74").Value = 15688527
$ws.Range("K74").Value = 15688527
$ws.Range("M74").Value = -15687653
$ws.Range("H77").Value = 12601015
$ws.Range("I77").Value = 15688527
$ws.Range("K77").Value = 78442635
$ws.Range("M77").Value = -78438267
$ws.Range("H136").Value = 125263450
$ws.Range("J136").Value = 550007
$ws.Range("L136").Value = 1650021
$ws.Range("N136").Value = -1655121

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3283.182
$ws.Range("I31").Value = 1301.579
$ws.Range("J31").Value = 15833.333
$ws.Range("K31").Value = 1301.579
$ws.Range("L31").Value = 15833.333
$ws.Range("M31").Value = -1006.579
$ws.Range("N31").Value = -16423.333
$ws.Range("H34").Value = 3283.182
$ws.Range("I34").Value = 1301.579
$ws.Range("J34").Value = 15833.333
$ws.Range("K34").Value = 1301.579
$ws.Range("L34").Value = 15833.333
$ws.Range("M34").Value = -1099.579
$ws.Range("N34").Value = -16237.333
$ws.Range("H58").Value = 23257434
$ws.Range("J58").Value = 2819.5334
$ws.Range("L58").Value = 2819.5334
$ws.Range("N58").Value = -3225.5334
$ws.Range("H134").Value = 47086.543
$ws.Range("I134").Value = 2262.2
$ws.Range("K134").Value = 6786.599999999999
$ws.Range("M134").Value = -4251.599999999999
$ws.Range("H136").Value = 23257434
$ws.Range("J136").Value = 2819.5334
$ws.Range("L136").Value = 8458.600199999999
$ws.Range("N136").Value = -13558.6002

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H24").Value = 0
$ws.Range("J24").Value = 0
$ws.Range("L24").Value = 0
$ws.Range("N24").ClearContents()
$ws.Range("H40").Value = 60.545456
$ws.Range("I40").Value = 60.545456
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 242.181824
$ws.Range("L40").Value = 0
$ws.Range("M40").Value = -173.181824
$ws.Range("N40").ClearContents()
$ws.Range("H86").Value = 1590.3846
$ws.Range("J86").Value = 1689.5834
$ws.Range("L86").Value = 5068.7502
$ws.Range("N86").Value = -7440.7502
$ws.Range("H89").Value = 1590.3846
$ws.Range("J89").Value = 1689.5834
$ws.Range("L89").Value = 15206.2506
$ws.Range("N89").Value = -27062.2506
$ws.Range("H103").Value = 2531.6
$ws.Range("I103").Value = 422.66666
$ws.Range("J103").Value = 4257.091
$ws.Range("K103").Value = 1267.99998
$ws.Range("L103").Value = 12771.273
$ws.Range("M103").Value = -388.9999800000001
$ws.Range("N103").Value = -14529.273
$ws.Range("H106").Value = 3350
$ws.Range("I106").Value = 3000
$ws.Range("J106").Value = 3466.6667
$ws.Range("K106").Value = 9000
$ws.Range("L106").Value = 10400.0001
$ws.Range("M106").Value = -8054
$ws.Range("N106").Value = -12292.0001
$ws.Range("H109").Value = 1616
$ws.Range("I109").Value = 1040
$ws.Range("J109").Value = 3200
$ws.Range("K109").Value = 3120
$ws.Range("L109").Value = 9600
$ws.Range("M109").Value = -2080
$ws.Range("N109").Value = -11680
$ws.Range("H112").Value = 13892211
$ws.Range("I112").Value = 1889.5714
$ws.Range("J112").Value = 19611756
$ws.Range("K112").Value = 5668.7142
$ws.Range("L112").Value = 58835268
$ws.Range("M112").Value = -4560.7142
$ws.Range("N112").Value = -58837484
$ws.Range("H115").Value = 2611.1428
$ws.Range("I115").Value = 2135.6
$ws.Range("J115").Value = 3800
$ws.Range("K115").Value = 6406.799999999999
$ws.Range("L115").Value = 11400
$ws.Range("M115").Value = -5231.799999999999
$ws.Range("N115").Value = -13750
$ws.Range("H118").Value = 3644.8276
$ws.Range("I118").Value = 1800
$ws.Range("J118").Value = 3710.7144
$ws.Range("K118").Value = 5400
$ws.Range("L118").Value = 11132.1432
$ws.Range("M118").Value = -4157
$ws.Range("N118").Value = -13618.1432
$ws.Range("H121").Value = 30735756
$ws.Range("I121").Value = 1067
$ws.Range("J121").Value = 32353372
$ws.Range("K121").Value = 3201
$ws.Range("L121").Value = 97060116
$ws.Range("M121").Value = -1891
$ws.Range("N121").Value = -97062736
$ws.Range("H122").Value = 717.04
$ws.Range("I122").Value = 322.07693
$ws.Range("J122").Value = 1144.9166
$ws.Range("K122").Value = 2898.69237
$ws.Range("L122").Value = 10304.2494
$ws.Range("M122").Value = -448.6923700000002
$ws.Range("N122").Value = -15204.2494
$ws.Range("H131").Value = 1304.6875
$ws.Range("I131").Value = 671
$ws.Range("J131").Value = 1592.7273
$ws.Range("K131").Value = 2013
$ws.Range("L131").Value = 4778.1819
$ws.Range("M131").Value = 3027
$ws.Range("N131").Value = -14858.1819

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 2384.077
$ws.Range("J126").Value = 2554.7778
$ws.Range("L126").Value = 7664.3334
$ws.Range("N126").Value = -12604.3334

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 114164.445
$ws.Range("I132").Value = 1875
$ws.Range("K132").Value = 5625
$ws.Range("M132").Value = -3095
$ws.Range("H136").Value = 106347.84
$ws.Range("I136").Value = 143742.72
$ws.Range("J136").Value = 84534.164
$ws.Range("K136").Value = 431228.16
$ws.Range("L136").Value = 253602.492
$ws.Range("M136").Value = -428678.16
$ws.Range("N136").Value = -258702.492
